# Add a new "31-jul" data column (AL) to the right of the existing date
# columns (A:AK), with a header in row 1 and one integer value per product
# row (2-11), mirroring the formatting already used by the "28-jul" (AK)
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: text-formatted like the other date headers in row 1.
$ws.Range("AL1").NumberFormat = "@"
$ws.Range("AL1").Value = "31-jul"

# Data cells: centered whole numbers, like AJ2:AK11.
$ws.Range("AL2:AL11").HorizontalAlignment = -4108
$ws.Range("AL2:AL11").NumberFormat = "0"

$ws.Range("AL2").Value = 13
$ws.Range("AL3").Value = 17
$ws.Range("AL4").Value = 9
$ws.Range("AL5").Value = 8
$ws.Range("AL6").Value = 12
$ws.Range("AL7").Value = 14
$ws.Range("AL8").Value = 11
$ws.Range("AL9").Value = 12
$ws.Range("AL10").Value = 17
$ws.Range("AL11").Value = 14

# Reposition the active selection as it was left in the saved workbook.
$ws.Range("AK16").Select() | Out-Null
